$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (new date string, new EBITDA string or $null if unchanged)
$updates = @{
    2  = @("2025/12/12", "5.80")
    8  = @("2025/12/12", "8.13")
    14 = @("2025/12/12", "2.94")
    20 = @("2025/12/12", "12.85")
    26 = @("2025/12/12", "10.81")
    32 = @("2025/12/12", "27.29")
    38 = @("2025/12/12", $null)
    44 = @("2025/12/12", "11.89")
    50 = @("2025/12/12", $null)
    56 = @("2025/12/12", "36.75")
    62 = @("2025/12/12", "12.64")
    68 = @("2025/12/12", "14.64")
    74 = @("2025/12/12", "16.96")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dateVal = $vals[0]
    $ebitdaVal = $vals[1]

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dateVal
    $cellA.Style = "Normal"

    if ($ebitdaVal -ne $null) {
        $cellB = $ws.Cells.Item($row, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $ebitdaVal
        $cellB.Style = "Normal"
    }
}
